$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'285.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-10.74%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.38%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07253"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-6.19%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.301"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.64%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.497"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-11.39%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9144"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.46%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-3.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1702"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-7.10%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08578"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04160"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.01%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-3.15%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005897"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.404"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.92%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E18").Value = "'-2.93%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.816"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.10%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1354"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.72%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2889"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.03848"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.58%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'0.44%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-8.16%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001285"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.12%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003733"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02284"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-10.02%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.04917"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-7.87%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.006699"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'236.35%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007705"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.11%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1265"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-4.00%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007379"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.006941"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-8.24%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3086"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-10.61%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006395"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.52%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.22%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'29.68%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.22%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.22%"
$ws.Range("E51").Style = "Normal"
